# Mark "Session 15" (column U) attendance as Absent ("A") for the
# students who were marked present by default but should be absent,
# per the final pre-publish attendance correction.
#
# Rows correspond to: 2024PGP005, 2021IPM013, 2021IPM018, 2021IPM022,
# 2024PGP095, 2024PGP506, 2024PGP154, 2024PGP164, 2024PGP166, 2024PGP177,
# 2024PGP196, 2021IPM063, 2024PGP560, 2024PGP465, 2024PGP258, 2024PGP531,
# 2024PGP282, 2024PGP558, 2024PGP292, 2021IPM103, 2024PGP295, 2024PGP336,
# 2024PGP349, 2024PGP356, 2024PGP519, 2024PGP378, 2024PGP394, 2024PGP452,
# 2024PGP540, 2024PGP462

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(8,13,19,20,21,24,27,28,29,32,34,35,39,41,44,45,46,47,49,50,51,55,57,58,60,61,66,71,72,73)

foreach ($r in $rows) {
    $cell = $ws.Range("U$r")
    $cell.Value = "A"
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
    $cell.WrapText = $true
}
